# Ajout d'une ligne dans le cahier des charges
# Adds a new row (14) to the "Feuil1" worksheet with a new
# "Fonction" / "Solution" / "Priorite" entry, matching the
# three new data cells introduced by the commit (A14, B14, D14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Accentuer l'immersion dans la partie"
$ws.Range("B14").Value = "Cacher un maximum les interfaces, préférer des raccourcis ou des actions simples à la souris"
$ws.Range("D14").Value = "Haute"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("B14").Select()
